# Update countries & provincias Spain
# Refresh the COVID-19 country data snapshot (new timestamp + updated
# case/recovered/death counts) and re-insert a few countries (Polonia,
# Uganda, San Martin) at their newly sorted position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner (row 1)
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 27 de Marzo de 2020 a las 08:42'

# Row 15
$ws.Cells.Item(15, 2).Value = 6962
$ws.Cells.Item(15, 3).Value = 53
$ws.Cells.Item(15, 5).Value = 6801

# Row 21
$ws.Cells.Item(21, 2).Value = 3166
$ws.Cells.Item(21, 3).Value = 116
$ws.Cells.Item(21, 5).Value = 2983

# Row 33
$ws.Cells.Item(33, 1).Value = 'Polonia'
$ws.Cells.Item(33, 2).Value = 1244
$ws.Cells.Item(33, 3).Value = 23
$ws.Cells.Item(33, 4).Value = 7
$ws.Cells.Item(33, 5).Value = 1221
$ws.Cells.Item(33, 6).Value = 3
$ws.Cells.Item(33, 8).Value = 16

# Row 34
$ws.Cells.Item(34, 1).Value = 'Pakistan'
$ws.Cells.Item(34, 2).Value = 1235
$ws.Cells.Item(34, 3).Value = 34
$ws.Cells.Item(34, 4).Value = 23
$ws.Cells.Item(34, 5).Value = 1203
$ws.Cells.Item(34, 6).Value = 7
$ws.Cells.Item(34, 8).Value = 9

# Row 66
$ws.Cells.Item(66, 2).Value = 345
$ws.Cells.Item(66, 3).Value = 46
$ws.Cells.Item(66, 5).Value = 340
$ws.Cells.Item(66, 6).Value = 2

# Row 86
$ws.Cells.Item(86, 5).Value = 150
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = 7

# Row 104
$ws.Cells.Item(104, 5).Value = 84
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = 2

# Row 112
$ws.Cells.Item(112, 2).Value = 70
$ws.Cells.Item(112, 3).Value = 1
$ws.Cells.Item(112, 5).Value = 69

# Row 140
$ws.Cells.Item(140, 1).Value = 'Uganda'
$ws.Cells.Item(140, 2).Value = 18
$ws.Cells.Item(140, 3).Value = 4
$ws.Cells.Item(140, 5).Value = 18

# Row 141
$ws.Cells.Item(141, 1).Value = 'Islas Virgenes de los Estados Unidos'
$ws.Cells.Item(141, 2).Value = 17
$ws.Cells.Item(141, 5).Value = 17

# Row 142
$ws.Cells.Item(142, 1).Value = 'Zambia'
$ws.Cells.Item(142, 2).Value = 16
$ws.Cells.Item(142, 4).Value = 0
$ws.Cells.Item(142, 5).Value = 16

# Row 143
$ws.Cells.Item(143, 1).Value = 'Bermudas'
$ws.Cells.Item(143, 2).Value = 15
$ws.Cells.Item(143, 4).Value = 2
$ws.Cells.Item(143, 5).Value = 13

# Row 147
$ws.Cells.Item(147, 4).Value = 9
$ws.Cells.Item(147, 5).Value = 4

# Row 151
$ws.Cells.Item(151, 1).Value = 'San Martin (Parte Francesa)'

# Row 153
$ws.Cells.Item(153, 1).Value = 'Republica de Yibuti'
